$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock List")

# Update "Last Updated" timestamp in Metadata sheet
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 02:45 PM"

# Stock List data shifts up by one row (oldest entry removed, new entry appended at bottom)
$stockData = @(
    @{B="NIFTYCASE"; C="NIFTYCASE"; D=10.19; E=-0.5854; H=0},
    @{B="MOMENTUM30"; C="MOMENTUM30"; D=31.54; E=-0.6614; H=0},
    @{B="CANHLIFE"; C="CANHLIFE"; D=118.46; E=0.6286; H=11253.7},
    @{B="FLEXIADD"; C="FLEXIADD"; D=10.64; E=-1.0233; H=0},
    @{B="MOENERGY"; C="MOENERGY"; D=36.3; E=-0.6568000000000001; H=0},
    @{B="MONIFTY100"; C="MONIFTY100"; D=26.49; E=0.3409; H=0},
    @{B="RUBICON"; C="RUBICON"; D=652.65; E=-0.1453; H=10752.4289},
    @{B="CRAMC"; C="CRAMC"; D=317.2; E=2.3226; H=6325.5208},
    @{B="LGEINDIA"; C="LGEINDIA"; D=1633.4; E=-0.946; H=110870.6825},
    @{B="TATACAP"; C="TATACAP"; D=329.3; E=0.1521; H=139783.5374},
    @{B="ELIQUID"; C="ELIQUID"; D=1004.85; E=0.0408; H=0},
    @{B="WEWORK"; C="WEWORK"; D=632.15; E=-2.4008; H=8472.2803},
    @{B="GROWWRLTY"; C="GROWWRLTY"; D=10.8; E=-0.4608; H=0},
    @{B="ADVANCE"; C="ADVANCE"; D=130.05; E=-5.2666; H=836.0358},
    @{B="OMFREIGHT"; C="OMFREIGHT"; D=88.90000000000001; E=-0.5926; H=299.3747},
    @{B="GLOTTIS"; C="GLOTTIS"; D=72.73999999999999; E=-0.8587; H=672.1394},
    @{B="FABTECH"; C="FABTECH"; D=237.72; E=0.4734; H=1056.6843},
    @{B="PACEDIGITK"; C="PACEDIGITK"; D=218.85; E=0.1327; H=4723.9063},
    @{B="JAINREC"; C="JAINREC"; D=377.25; E=1.2208; H=13018.3623},
    @{B="EPACKPEB"; C="EPACKPEB"; D=301.45; E=1.979; H=3028.1254},
    @{B="BMWVENTLTD"; C="BMWVENTLTD"; D=69.25; E=0; H=600.5014},
    @{B="STYL"; C="STYL"; D=372.4; E=-0.8388; H=6025.649},
    @{B="JARO"; C="JARO"; D=621.5; E=-1.4821; H=1377.0134},
    @{B="SOLARWORLD"; C="SOLARWORLD"; D=309.1; E=-0.6269; H=2679.0517},
    @{B="ARSSBL"; C="ARSSBL"; D=537.3; E=4.7266; H=3370.2277},
    @{B="GANESHCP"; C="GANESHCP"; D=274.4; E=-2.7984; H=1108.9312},
    @{B="ATLANTAELE"; C="ATLANTAELE"; D=1003.05; E=-1.7436; H=7713.116},
    @{B="GKENERGY"; C="GKENERGY"; D=213.85; E=-0.7933; H=4337.2472},
    @{B="SAATVIKGL"; C="SAATVIKGL"; D=528.2; E=-1.3079; H=6713.6863},
    @{B="IVALUE"; C="IVALUE"; D=281.45; E=-0.3364; H=1506.8799},
    @{B="VMSTMT"; C="VMSTMT"; D=70.03; E=-0.9056; H=347.5674},
    @{B="EUROPRATIK"; C="EUROPRATIK"; D=321.75; E=0.8147; H=3288.285},
    @{B="SHRINGARMS"; C="SHRINGARMS"; D=229.31; E=-1.2616; H=2211.284},
    @{B="DEVX"; C="DEVX"; D=44.53; E=-0.3803; H=401.605},
    @{B="URBANCO"; C="URBANCO"; D=148.9; E=-2.0459; H=21380.5798},
    @{B="SML100CASE"; C="SML100CASE"; D=10.36; E=-0.7663; H=0},
    @{B="AONEGOLD"; C="AONEGOLD"; D=11.28; E=-0.2653; H=0},
    @{B="ELM250"; C="ELM250"; D=16.72; E=0.1797; H=0},
    @{B="AMANTA"; C="AMANTA"; D=122.52; E=1.407; H=475.7372},
    @{B="CPEDU"; C="CPEDU"; D=315.9; E=1.8539; H=574.7148999999999},
    @{B="AHCL"; C="AHCL"; D=139.27; E=3.1706; H=740.2409},
    @{B="STLNETWORK"; C="STLNETWORK"; D=26.59; E=-0.412; H=1297.3822},
    @{B="VIKRAN"; C="VIKRAN"; D=98.05; E=-1.783; H=2528.8166},
    @{B="MANUFGBEES"; C="MANUFGBEES"; D=151.77; E=-1.011; H=0},
    @{B="MEIL"; C="MEIL"; D=461.15; E=-0.7319; H=1274.1632},
    @{B="GROWWNXT50"; C="GROWWNXT50"; D=70.29000000000001; E=-0.4109; H=0},
    @{B="SHREEJISPG"; C="SHREEJISPG"; D=270.05; E=-0.7899; H=4399.6074},
    @{B="GEMAROMA"; C="GEMAROMA"; D=219.52; E=-0.876; H=1146.7097},
    @{B="PATELRMART"; C="PATELRMART"; D=219.31; E=-1.0646; H=732.5069999999999},
    @{B="VIKRAMSOLR"; C="VIKRAMSOLR"; D=322; E=-1.5892; H=11647.2884},
    @{B="LTGILTCASE"; C="LTGILTCASE"; D=29.67; E=0.2365; H=0},
    @{B="REGAAL"; C="REGAAL"; D=89.13; E=-0.8675; H=915.5742},
    @{B="BLUESTONE"; C="BLUESTONE"; D=711.95; E=0.1266; H=10773.2539},
    @{B="MOSILVER"; C="MOSILVER"; D=145.9; E=-1.5054; H=0},
    @{B="ALLTIME"; C="ALLTIME"; D=308.75; E=2.66; H=2022.5526},
    @{B="JSWCEMENT"; C="JSWCEMENT"; D=134.98; E=-0.4793; H=18402.6999},
    @{B="SBILIQETF"; C="SBILIQETF"; D=1012.94; E=0.0296; H=0},
    @{B="HILINFRA"; C="HILINFRA"; D=77.23; E=-0.3998; H=0},
    @{B="GROWWPOWER"; C="GROWWPOWER"; D=10.28; E=-0.9634; H=0},
    @{B="LOTUSDEV"; C="LOTUSDEV"; D=177.82; E=0.3669; H=8690.485000000001},
    @{B="MBEL"; C="MBEL"; D=450.2; E=-0.7714; H=2572.8126},
    @{B="LAXMIINDIA"; C="LAXMIINDIA"; D=145.62; E=-1.1942; H=761.1248000000001},
    @{B="CPPLUS"; C="CPPLUS"; D=1322.1; E=-0.264; H=15497.9053},
    @{B="SHANTIGOLD"; C="SHANTIGOLD"; D=241.57; E=-1.6409; H=1741.6231},
    @{B="MOGOLD"; C="MOGOLD"; D=119.65; E=-0.5403; H=0},
    @{B="BRIGHOTEL"; C="BRIGHOTEL"; D=82.39; E=-0.9855; H=3129.5229},
    @{B="INDIQUBE"; C="INDIQUBE"; D=212.64; E=-0.7561; H=4465.6847},
    @{B="EBGNG"; C="EBGNG"; D=346.65; E=3.2311; H=3952.2092},
    @{B="LIQGRWBEES"; C="LIQGRWBEES"; D=1014.74; E=0.0246; H=0},
    @{B="CHEMBONDCH"; C="CHEMBONDCH"; D=153.35; E=-1.6987; H=412.459},
    @{B="GROWWNIFTY"; C="GROWWNIFTY"; D=10.29; E=-0.3872; H=0},
    @{B="ANTHEM"; C="ANTHEM"; D=702.25; E=-0.1209; H=39439.0658},
    @{B="QUALITY30"; C="QUALITY30"; D=21.05; E=-0.8945; H=0},
    @{B="SMARTWORKS"; C="SMARTWORKS"; D=606.65; E=2.0867; H=6931.2448},
    @{B="TRAVELFOOD"; C="TRAVELFOOD"; D=1316.3; E=0.1141; H=17332.9705}
)

$r = 2
foreach ($item in $stockData) {
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 8).Value = $item.H
    $r++
}